$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '42.033.35'
$ws.Range('E2').Value = '  -0.91%  '

# Row 3
$ws.Range('D3').Value = '2.218.76'
$ws.Range('E3').Value = '  -1.63%  '

# Row 4
$ws.Range('E4').Value = '  +0.09%  '

# Row 5
$ws.Range('E5').Value = '  -2.03%  '

# Row 6
$ws.Range('E6').Value = '  +0.81%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '74.00'
$ws.Range('E7').Value = '  -0.91%  '

# Row 8
$ws.Range('E8').Value = '  +0.06%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.619'
$ws.Range('E9').Value = '  -0.57%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '43.75'
$ws.Range('E10').Value = '  +5.08%  '

# Row 11
$ws.Range('E11').Value = '  +2.07%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.12'
$ws.Range('E12').Value = '  -0.17%  '

# Row 13
$ws.Range('E13').Value = '  +0.16%  '

# Row 14
$ws.Range('D14').Value = '2.550.66'
$ws.Range('E14').Value = '  -1.50%  '

# Row 15
$ws.Range('E15').Value = '  -1.30%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.25'
$ws.Range('E16').Value = '  -1.72%  '

# Row 17
$ws.Range('D17').Value = '2.214.50'
$ws.Range('E17').Value = '  -1.75%  '

# Row 18
$ws.Range('D18').Value = '41.855.68'
$ws.Range('E18').Value = '  -0.91%  '

# Row 19
$ws.Range('E19').Value = '  +12.01%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.20'
$ws.Range('E20').Value = '  +1.06%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.40'
$ws.Range('E21').Value = '  +0.70%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.49'
$ws.Range('E22').Value = '  +31.84%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '229.76'
$ws.Range('E23').Value = '  -0.27%  '

# Row 24
$ws.Range('E24').Value = '  -6.81%  '

# Row 25
$ws.Range('E25').Value = '  +0.01%  '

# Row 26
$ws.Range('E26').Value = '  +2.98%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.59'
$ws.Range('E27').Value = '  +1.17%  '

# Row 28
$ws.Range('E28').Value = '  -1.86%  '

# Row 29
$ws.Range('E29').Value = '  -0.16%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '166.61'
$ws.Range('E30').Value = '  -1.54%  '

# Row 31
$ws.Range('E31').Value = '  -0.68%  '

# Row 32
$ws.Range('E32').Value = '  +16.16%  '

# Row 33
$ws.Range('E33').Value = '  -3.21%  '

# Row 34
$ws.Range('E34').Value = '  -0.58%  '

# Row 35
$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '29.37'
$ws.Range('E35').Value = '  -3.84%  '

# Row 36
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.114'
$ws.Range('E36').Value = '  -4.60%  '

# Row 37
$ws.Range('E37').Value = '  -5.03%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0301'
$ws.Range('E38').Value = '  -0.90%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '13.01'
$ws.Range('E39').Value = '  -4.84%  '

# Row 40
$ws.Range('E40').Value = '  -2.04%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '65.34'
$ws.Range('E41').Value = '  +5.84%  '

# Row 42
$ws.Range('E42').Value = '  -2.25%  '

# Row 43
$ws.Range('E43').Value = '  -1.95%  '

# Row 44
$ws.Range('E44').Value = '  +0.65%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '104.37'
$ws.Range('E45').Value = '  -3.91%  '

# Row 46
$ws.Range('E46').Value = '  +0.22%  '

# Row 47
$ws.Range('E47').Value = '  +6.04%  '

# Row 48
$ws.Range('E48').Value = '  -0.67%  '

# Row 49
$ws.Range('E49').Value = '  -0.40%  '

# Row 50
$ws.Range('E50').Value = '  +0.57%  '

# Row 51
$ws.Range('D51').Value = '2.426.95'
$ws.Range('E51').Value = '  -1.48%  '
